$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC3").Value = 81
$ws.Range("AI3").Value = 34
$ws.Range("G3").Value = 2.4
$ws.Range("I3").Value = 3.1
$ws.Range("N3").Value = 2.75
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 2.2
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 6
$ws.Range("U3").Value = 10
$ws.Range("AD4").Value = 501
$ws.Range("AJ4").Value = 51
$ws.Range("H4").Value = 3.9
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 1.07
$ws.Range("K4").Value = 9
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 1.67
$ws.Range("Y4").Value = 34
$ws.Range("Z4").Value = 9
$ws.Range("AA5").Value = 7.5
$ws.Range("AB5").Value = 15
$ws.Range("AD5").Value = 201
$ws.Range("AF5").Value = 23
$ws.Range("AG5").Value = 15
$ws.Range("AI5").Value = 34
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4.5
$ws.Range("P5").Value = 1.33
$ws.Range("Q5").Value = 3.25
$ws.Range("R5").Value = 1.73
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 8
$ws.Range("U5").Value = 8.5
$ws.Range("W5").Value = 13
$ws.Range("Y5").Value = 23
$ws.Range("Z5").Value = 13
$ws.Range("AA6").Value = 8.5
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 40
$ws.Range("AD6").Value = 200
$ws.Range("AE6").Value = 19.5
$ws.Range("AF6").Value = 37
$ws.Range("AG6").Value = 16.5
$ws.Range("AJ6").Value = 35
$ws.Range("G6").Value = 1.37
$ws.Range("H6").Value = 4.55
$ws.Range("I6").Value = 6.2
$ws.Range("T6").Value = 8.25
$ws.Range("U6").Value = 7.1
$ws.Range("V6").Value = 7.2
$ws.Range("W6").Value = 8.5
$ws.Range("X6").Value = 8.75
$ws.Range("Y6").Value = 16
$ws.Range("Z6").Value = 17.5
$ws.Range("R7").Value = 1.89
$ws.Range("S7").Value = 1.82
$ws.Range("AB8").Value = 11
$ws.Range("AC8").Value = 34
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 19
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 34
$ws.Range("AI8").Value = 23
$ws.Range("AJ8").Value = 26
$ws.Range("G8").Value = 2.1
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 1.03
$ws.Range("K8").Value = 15
$ws.Range("L8").Value = 1.18
$ws.Range("M8").Value = 4.5
$ws.Range("P8").Value = 1.3
$ws.Range("Q8").Value = 3.4
$ws.Range("R8").Value = 1.53
$ws.Range("S8").Value = 2.38
$ws.Range("U8").Value = 12
$ws.Range("W8").Value = 21
$ws.Range("AA9").Value = 7
$ws.Range("AE9").Value = 10
$ws.Range("AF9").Value = 17
$ws.Range("AG9").Value = 11
$ws.Range("AH9").Value = 34
$ws.Range("AI9").Value = 23
$ws.Range("AJ9").Value = 29
$ws.Range("G9").Value = 2.05
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 1.04
$ws.Range("K9").Value = 13
$ws.Range("M9").Value = 3.75
$ws.Range("N9").Value = 1.88
$ws.Range("O9").Value = 1.98
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 1.73
$ws.Range("S9").Value = 2
$ws.Range("U9").Value = 11
$ws.Range("V9").Value = 9
$ws.Range("W9").Value = 19
$ws.Range("X9").Value = 17
$ws.Range("Y9").Value = 26
$ws.Range("Z9").Value = 12
$ws.Range("AA10").Value = 8.5
$ws.Range("AF10").Value = 26
$ws.Range("AG10").Value = 15
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 3.9
$ws.Range("J10").Value = 1.02
$ws.Range("K10").Value = 19
$ws.Range("L10").Value = 1.14
$ws.Range("M10").Value = 5.5
$ws.Range("N10").Value = 1.5
$ws.Range("O10").Value = 2.5
$ws.Range("P10").Value = 1.25
$ws.Range("Q10").Value = 3.75
$ws.Range("R10").Value = 1.53
$ws.Range("S10").Value = 2.38
$ws.Range("T10").Value = 11
$ws.Range("X10").Value = 12
$ws.Range("Y10").Value = 19
$ws.Range("Z10").Value = 19
$ws.Range("J11").Value = 1.04
$ws.Range("K11").Value = 13
$ws.Range("L11").Value = 1.22
$ws.Range("M11").Value = 4
$ws.Range("P11").Value = 1.33
$ws.Range("Q11").Value = 3.25
$ws.Range("R11").Value = 1.67
$ws.Range("S11").Value = 2.1
$ws.Range("U11").Value = 9
$ws.Range("J12").Value = 1.06
$ws.Range("K12").Value = 10
$ws.Range("L12").Value = 1.29
$ws.Range("M12").Value = 3.5
$ws.Range("P12").Value = 1.4
$ws.Range("Q12").Value = 2.75
$ws.Range("R12").Value = 1.73
$ws.Range("V12").Value = 9
$ws.Range("X12").Value = 17
$ws.Range("Y12").Value = 26
$ws.Range("AA13").Value = 6.1
$ws.Range("AD13").Value = 700
$ws.Range("AI13").Value = 30
$ws.Range("H13").Value = 3.1
$ws.Range("R13").Value = 1.87
$ws.Range("S13").Value = 1.83
$ws.Range("T13").Value = 6.9
$ws.Range("V13").Value = 9
$ws.Range("Y13").Value = 32
$ws.Range("G14").Value = 2.15
$ws.Range("H14").Value = 3.3
$ws.Range("K14").Value = 9.5
$ws.Range("L14").Value = 1.3
$ws.Range("M14").Value = 3.4
$ws.Range("N14").Value = 2.05
$ws.Range("O14").Value = 1.72
$ws.Range("P14").Value = 1.44
$ws.Range("Q14").Value = 2.63
$ws.Range("Y14").Value = 29
$ws.Range("Z14").Value = 9.5
$ws.Range("R15").Value = 1.73
$ws.Range("S15").Value = 1.99
$ws.Range("P17").Value = 1.47
$ws.Range("N18").Value = 1.48
$ws.Range("O18").Value = 2.6
$ws.Range("P18").Value = 1.22
$ws.Range("K20").Value = 10
$ws.Range("P20").Value = 1.37
$ws.Range("AA22").Value = 8
$ws.Range("AB22").Value = 19.5
$ws.Range("AC22").Value = 100
$ws.Range("AD22").Value = 900
$ws.Range("AE22").Value = 16
$ws.Range("AF22").Value = 40
$ws.Range("AG22").Value = 21
$ws.Range("AH22").Value = 150
$ws.Range("AI22").Value = 80
$ws.Range("AJ22").Value = 75
$ws.Range("G22").Value = 1.45
$ws.Range("H22").Value = 4.1
$ws.Range("I22").Value = 6.7
$ws.Range("J22").Value = 1.05
$ws.Range("K22").Value = 7.6
$ws.Range("L22").Value = 1.27
$ws.Range("M22").Value = 3.4
$ws.Range("N22").Value = 1.82
$ws.Range("O22").Value = 1.9
$ws.Range("P22").Value = 1.4
$ws.Range("Q22").Value = 2.75
$ws.Range("R22").Value = 1.98
$ws.Range("S22").Value = 1.75
$ws.Range("U22").Value = 6.5
$ws.Range("V22").Value = 8.25
$ws.Range("W22").Value = 9.5
$ws.Range("X22").Value = 12
$ws.Range("Y22").Value = 29
$ws.Range("Z22").Value = 7.6
$ws.Range("AA23").Value = 7
$ws.Range("AB23").Value = 17
$ws.Range("AD23").Value = 351
$ws.Range("AE23").Value = 12
$ws.Range("AF23").Value = 23
$ws.Range("AG23").Value = 15
$ws.Range("AH23").Value = 51
$ws.Range("AI23").Value = 41
$ws.Range("G23").Value = 1.75
$ws.Range("H23").Value = 3.6
$ws.Range("I23").Value = 4.75
$ws.Range("J23").Value = 1.06
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = 1.33
$ws.Range("M23").Value = 3.25
$ws.Range("N23").Value = 2.05
$ws.Range("O23").Value = 1.75
$ws.Range("P23").Value = 1.44
$ws.Range("Q23").Value = 2.63
$ws.Range("R23").Value = 1.91
$ws.Range("S23").Value = 1.8
$ws.Range("T23").Value = 6.5
$ws.Range("U23").Value = 8
$ws.Range("W23").Value = 13
$ws.Range("Y23").Value = 29
$ws.Range("Z23").Value = 9.5
